$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author selected column B's data range (B2:B101) before making the edit,
# which is why the saved sheetView shows that as the active selection.
$ws.Range("B2:B101").Select()

# The edit itself: every image path stored in column B gets a "..\" prefix
# inserted immediately before "static\" (static\К1.jpg -> ..\static\К1.jpg,
# and so on through К101.jpg). Using Range.Replace (rather than per-cell
# Value assignment) preserves the shared-string table's existing order/ids,
# matching the original author's commit which rewrote the <si> entries in
# place instead of appending new ones.
$ws.Range("B2:B101").Replace("static\", "..\static\")
